$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G width (closest achievable snap to ~10.42578125 chars) ---
$ws.Columns("G").ColumnWidth = 9.7

# --- F column (Material cost) actuals ---
$ws.Range("F5").Value = 7000
$ws.Range("F7").Value = 2000
$ws.Range("F8").Value = 10000
$ws.Range("F10").Value = 47000
$ws.Range("F11").Value = 3000
$ws.Range("F13").Value = 14000
$ws.Range("F14").Value = 1550
$ws.Range("F15").Value = 5190
$ws.Range("F16").Value = 2500
$ws.Range("F17").Value = 2200
$ws.Range("F18").Value = 25000
$ws.Range("F19").Value = 800
$ws.Range("F20").Value = 12000
$ws.Range("F21").Value = 3400
$ws.Range("F22").Value = 1750
$ws.Range("F23").Value = 3500
$ws.Range("F24").Value = 7900
$ws.Range("F25").Value = 2400
$ws.Range("F26").Value = 2400
$ws.Range("F27").Value = 12000

# --- Column G (variance = D - F) formulas ---
$ws.Range("G5").Formula = "=D5-F5"
$ws.Range("G6").Formula = "=D6-F6"
$ws.Range("G7").Formula = "=D7-F7"
$ws.Range("G8").Formula = "=D8-F8"
$ws.Range("G10").Formula = "=D10-F10"
$ws.Range("G13").Formula = "=D13-F13"
$ws.Range("G14").Formula = "=D14-F14"
$ws.Range("G15").Formula = "=D15-F15"
$ws.Range("G16").Formula = "=D16-F16"
$ws.Range("G17").Formula = "=D17-F17"
$ws.Range("G18").Formula = "=D18-F18"
$ws.Range("G19").Formula = "=D19-F19"
$ws.Range("G20").Formula = "=D20-F20"
$ws.Range("G21").Formula = "=D21-F21"
$ws.Range("G22").Formula = "=D22-F22"
$ws.Range("G23").Formula = "=D23-F23"
$ws.Range("G24").Formula = "=D24-F24"
$ws.Range("G25").Formula = "=D25-F25"
$ws.Range("G26").Formula = "=D26-F26"
$ws.Range("G27").Formula = "=D27-F27"

# G11 is formatted like the others but intentionally left blank (no formula/value)
$ws.Range("G11").NumberFormat = '_("$"* #,##0_);_("$"* \(#,##0\);_("$"* "-"??_);_(@_)'

# --- Highlight two line-item labels in yellow ---
$ws.Range("B14").Interior.Color = 65535
$ws.Range("B22").Interior.Color = 65535

# --- View state: scroll down and leave selection on G26 ---
$ws.Range("G26").Select()
